# Registracija Objekta.xlsx - apply commit changes
#
# Summary of the edit (per the xml diff):
#  1. Cell A11 text is shortened from
#     "2. Popunjava licne podatke i podatke u vezi objekta" to
#     "2. Popunjava licne podatke " (trailing space kept). Excel re-wraps
#     the cell and shrinks the row height back to the default 15pt.
#  2. The active pane's selection moves from B24 to G9 (and the view
#     scrolls down a row, topLeftCell A13 -> A14 - window-scroll position
#     is not something this content edit needs to drive explicitly, it is
#     a side effect of the user's navigation captured by Excel).
#  (All the other shared-string index churn in the diff is an automatic
#  side effect of Excel rebuilding the shared-strings table once the A11
#  string's content changes - no other cell's displayed text actually
#  changes.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shorten the "Popunjava licne podatke" step description.
$ws.Range("A11").Value = "2. Popunjava licne podatke "

# Excel re-fits the row to the (now single-line) text instead of the old
# two-line wrap, so the custom height drops back to the sheet default.
$ws.Rows.Item(11).RowHeight = 15

# 2) Move the selection/active cell to G9, as captured in the saved view.
$ws.Activate()
$ws.Range("G9").Select()
